$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 25 (old rows 25:31 shift down to 28:34)
$ws.Rows("25:27").Insert()

# New row 25
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "Femacal de La Calera"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 45142
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 100112043
$ws.Range("G25").Value = "Pepino dulce"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 56
$ws.Range("K25").Value = 23000
$ws.Range("L25").Value = 23000
$ws.Range("M25").Value = 23000
$ws.Range("N25").Value = "`$/caja 15 kilos"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 1533
$ws.Range("Q25").Value = 15
$ws.Range("R25").Value = "Hortaliza"

# New row 26
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 45142
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 100112043
$ws.Range("G26").Value = "Pepino dulce"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 21000
$ws.Range("L26").Value = 21000
$ws.Range("M26").Value = 21000
$ws.Range("N26").Value = "`$/caja 15 kilos"
$ws.Range("O26").Value = "Provincia de Petorca"
$ws.Range("P26").Value = 1400
$ws.Range("Q26").Value = 15
$ws.Range("R26").Value = "Hortaliza"

# New row 27
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 45142
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112043
$ws.Range("G27").Value = "Pepino dulce"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("N27").Value = "`$/caja 15 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 1200
$ws.Range("Q27").Value = 15
$ws.Range("R27").Value = "Hortaliza"
